$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3476551.2
$ws.Range("I15").Value = 3476551.2
$ws.Range("K15").Value = 10429653.6
$ws.Range("M15").Value = -10429484.6
$ws.Range("H40").Value = 4078.2144
$ws.Range("I40").Value = 3065.8333
$ws.Range("J40").Value = 4837.5
$ws.Range("K40").Value = 3065.8333
$ws.Range("L40").Value = 4837.5
$ws.Range("M40").Value = -2890.8333
$ws.Range("N40").Value = -5187.5
$ws.Range("H80").Value = 1407
$ws.Range("I80").Value = 665
$ws.Range("J80").Value = 2149
$ws.Range("K80").Value = 1995
$ws.Range("L80").Value = 6447
$ws.Range("M80").Value = -997
$ws.Range("N80").Value = -8443
$ws.Range("H83").Value = 1407
$ws.Range("I83").Value = 665
$ws.Range("J83").Value = 2149
$ws.Range("K83").Value = 5985
$ws.Range("L83").Value = 19341
$ws.Range("M83").Value = -993
$ws.Range("N83").Value = -29325
$ws.Range("H131").Value = 3796.4
$ws.Range("I131").Value = 1694.6
$ws.Range("J131").Value = 8000
$ws.Range("K131").Value = 5083.799999999999
$ws.Range("L131").Value = 24000
$ws.Range("M131").Value = -43.79999999999927
$ws.Range("N131").Value = -34080
$ws.Range("H132").Value = 1050.3823
$ws.Range("I132").Value = 1076.4667
$ws.Range("J132").Value = 854.75
$ws.Range("K132").Value = 3229.4001
$ws.Range("L132").Value = 2564.25
$ws.Range("M132").Value = -699.4000999999998
$ws.Range("N132").Value = -7624.25
$ws.Range("H137").Value = 6003.7144
$ws.Range("I137").Value = 3884.111
$ws.Range("J137").Value = 8248
$ws.Range("K137").Value = 11652.333
$ws.Range("L137").Value = 24744
$ws.Range("M137").Value = -9102.332999999999
$ws.Range("N137").Value = -29844
$ws.Range("H138").Value = 2829.3667
$ws.Range("I138").Value = 2552.2
$ws.Range("J138").Value = 2854.5637
$ws.Range("K138").Value = 7656.599999999999
$ws.Range("L138").Value = 8563.6911
$ws.Range("M138").Value = -2516.599999999999
$ws.Range("N138").Value = -18843.6911

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 1497.25
$ws.Range("I14").Value = 3399
$ws.Range("J14").Value = 356.2
$ws.Range("K14").Value = 3399
$ws.Range("L14").Value = 356.2
$ws.Range("M14").Value = -3224
$ws.Range("N14").Value = -706.2
$ws.Range("H32").Value = 8936316
$ws.Range("I32").Value = 10006216
$ws.Range("K32").Value = 10006216
$ws.Range("M32").Value = -10005929
$ws.Range("H61").Value = 20840126
$ws.Range("I61").Value = 15631066
$ws.Range("J61").Value = 62512612
$ws.Range("K61").Value = 15631066
$ws.Range("L61").Value = 62512612
$ws.Range("M61").Value = -15630854
$ws.Range("N61").Value = -62513036
$ws.Range("H136").Value = 20840126
$ws.Range("I136").Value = 15631066
$ws.Range("J136").Value = 62512612
$ws.Range("K136").Value = 46893198
$ws.Range("L136").Value = 187537836
$ws.Range("M136").Value = -46890648
$ws.Range("N136").Value = -187542936

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 44645.25
$ws.Range("J140").Value = 44645.25
$ws.Range("L140").Value = 44645.25
$ws.Range("N140").Value = -55005.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 786406.2
$ws.Range("I31").Value = 10453.883
$ws.Range("J31").Value = 1801113
$ws.Range("K31").Value = 10453.883
$ws.Range("L31").Value = 1801113
$ws.Range("M31").Value = -10158.883
$ws.Range("N31").Value = -1801703
$ws.Range("H34").Value = 786406.2
$ws.Range("I34").Value = 10453.883
$ws.Range("J34").Value = 1801113
$ws.Range("K34").Value = 10453.883
$ws.Range("L34").Value = 1801113
$ws.Range("M34").Value = -10251.883
$ws.Range("N34").Value = -1801517
$ws.Range("H35").Value = 2716.6667
$ws.Range("I35").Value = 2716.6667
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 2716.6667
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -2422.6667
$ws.Range("N35").ClearContents()
$ws.Range("H122").Value = 1160.5333
$ws.Range("I122").Value = 1200.5714
$ws.Range("J122").Value = 600
$ws.Range("K122").Value = 3601.7142
$ws.Range("L122").Value = 1800
$ws.Range("M122").Value = -1151.7142
$ws.Range("N122").Value = -6700
$ws.Range("H132").Value = 1848.875
$ws.Range("I132").Value = 1840.4193
$ws.Range("J132").Value = 2111
$ws.Range("K132").Value = 5521.257900000001
$ws.Range("L132").Value = 6333
$ws.Range("M132").Value = -2991.257900000001
$ws.Range("N132").Value = -11393

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1264
$ws.Range("I113").Value = 512.4
$ws.Range("J113").Value = 1577.1666
$ws.Range("K113").Value = 1537.2
$ws.Range("L113").Value = 4731.4998
$ws.Range("M113").Value = 632.8000000000002
$ws.Range("N113").Value = -9071.4998
$ws.Range("H137").Value = 4683.75
$ws.Range("I137").Value = 2550.6
$ws.Range("K137").Value = 7651.799999999999
$ws.Range("M137").Value = -2551.799999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2370.7896
$ws.Range("I122").Value = 2363.611
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 7090.833
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -4640.833
$ws.Range("N122").Value = -12400
$ws.Range("H136").Value = 33999.5
$ws.Range("J136").Value = 33999.5
$ws.Range("L136").Value = 101998.5
$ws.Range("N136").Value = -107098.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 44066.73
$ws.Range("I7").Value = 4823.1665
$ws.Range("J7").Value = 132364.75
$ws.Range("K7").Value = 4823.1665
$ws.Range("L7").Value = 132364.75
$ws.Range("M7").Value = -4711.1665
$ws.Range("N7").Value = -132588.75
$ws.Range("H22").Value = 1397.8334
$ws.Range("I22").Value = 1350
$ws.Range("J22").Value = 1421.75
$ws.Range("K22").Value = 1350
$ws.Range("L22").Value = 1421.75
$ws.Range("M22").Value = -1055
$ws.Range("N22").Value = -2011.75
$ws.Range("H27").Value = 1397.8334
$ws.Range("I27").Value = 1350
$ws.Range("J27").Value = 1421.75
$ws.Range("K27").Value = 1350
$ws.Range("L27").Value = 1421.75
$ws.Range("M27").Value = -1243
$ws.Range("N27").Value = -1635.75
$ws.Range("H36").Value = 99424
$ws.Range("J36").Value = 99424
$ws.Range("L36").Value = 99424
$ws.Range("N36").Value = -100548
$ws.Range("H43").Value = 1086573.4
$ws.Range("I43").Value = 38400
$ws.Range("J43").Value = 1324794.5
$ws.Range("K43").Value = 38400
$ws.Range("L43").Value = 1324794.5
$ws.Range("M43").Value = -38207
$ws.Range("N43").Value = -1325180.5
$ws.Range("H46").Value = 4517.185
$ws.Range("I46").Value = 1868.1
$ws.Range("J46").Value = 12086
$ws.Range("K46").Value = 1868.1
$ws.Range("L46").Value = 12086
$ws.Range("M46").Value = -1680.1
$ws.Range("N46").Value = -12462
$ws.Range("H122").Value = 5932.9443
$ws.Range("I122").Value = 5499.5835
$ws.Range("J122").Value = 6799.6665
$ws.Range("K122").Value = 16498.7505
$ws.Range("L122").Value = 20398.9995
$ws.Range("M122").Value = -14048.7505
$ws.Range("N122").Value = -25298.9995
$ws.Range("H126").Value = 44066.73
$ws.Range("I126").Value = 4823.1665
$ws.Range("J126").Value = 132364.75
$ws.Range("K126").Value = 14469.4995
$ws.Range("L126").Value = 397094.25
$ws.Range("M126").Value = -11999.4995
$ws.Range("N126").Value = -402034.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 104760
$ws.Range("J109").Value = 104760
$ws.Range("L109").Value = 104760
$ws.Range("N109").Value = -107534
$ws.Range("H122").Value = 1426.5217
$ws.Range("I122").Value = 1426.5217
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4279.5651
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1829.5651
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 2894
$ws.Range("I126").Value = 2641.1667
$ws.Range("J126").Value = 3399.6667
$ws.Range("K126").Value = 7923.500100000001
$ws.Range("L126").Value = 10199.0001
$ws.Range("M126").Value = -5453.500100000001
$ws.Range("N126").Value = -15139.0001
$ws.Range("H136").Value = 2980.4
$ws.Range("I136").Value = 2904.9644
$ws.Range("J136").Value = 3282.1428
$ws.Range("K136").Value = 8714.893199999999
$ws.Range("L136").Value = 9846.428400000001
$ws.Range("M136").Value = -6164.893199999999
$ws.Range("N136").Value = -14946.4284
